$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2169811320754717
$ws.Range("C2").Value = 0.5150943396226415
$ws.Range("J2").Value = 0.0169811320754717
$ws.Range("O2").Value = 0.001886792452830189
$ws.Range("P2").Value = 0.1471698113207547
$ws.Range("S2").Value = 0.1018867924528302

# Row 3
$ws.Range("B3").Value = 0.0176056338028169
$ws.Range("C3").Value = 0.04577464788732395
$ws.Range("J3").Value = 0.01056338028169014
$ws.Range("P3").Value = 0.721830985915493
$ws.Range("S3").Value = 0.2042253521126761

# Row 4
$ws.Range("J4").Value = 0.02272727272727273
$ws.Range("P4").Value = 0.7045454545454546
$ws.Range("S4").Value = 0.2727272727272727

# Row 6
$ws.Range("B6").Value = 0.06183368869936034
$ws.Range("D6").Value = 0.01066098081023454
$ws.Range("E6").Value = 0.002132196162046908
$ws.Range("F6").Value = 0.09594882729211088
$ws.Range("J6").Value = 0.2238805970149254
$ws.Range("O6").Value = 0.02132196162046908
$ws.Range("Q6").Value = 0.1513859275053305
$ws.Range("R6").Value = 0.05330490405117271
$ws.Range("S6").Value = 0.3795309168443497

# Row 7
$ws.Range("B7").Value = 0.08490566037735849
$ws.Range("D7").Value = 0.02830188679245283
$ws.Range("E7").Value = 0.002358490566037736
$ws.Range("F7").Value = 0.09433962264150944
$ws.Range("J7").Value = 0.1155660377358491
$ws.Range("O7").Value = 0.02358490566037736
$ws.Range("Q7").Value = 0.1863207547169811
$ws.Range("R7").Value = 0.08018867924528301
$ws.Range("S7").Value = 0.3844339622641509

# Row 8
$ws.Range("B8").Value = 0.07782101167315175
$ws.Range("D8").Value = 0.01750972762645914
$ws.Range("F8").Value = 0.0632295719844358
$ws.Range("J8").Value = 0.1021400778210117
$ws.Range("O8").Value = 0.0301556420233463
$ws.Range("Q8").Value = 0.1692607003891051
$ws.Range("R8").Value = 0.08754863813229571
$ws.Range("S8").Value = 0.4523346303501946

# Row 9
$ws.Range("B9").Value = 0.08624229979466119
$ws.Range("D9").Value = 0.01848049281314168
$ws.Range("F9").Value = 0.06365503080082136
$ws.Range("J9").Value = 0.08829568788501027
$ws.Range("O9").Value = 0.03490759753593429
$ws.Range("Q9").Value = 0.1765913757700205
$ws.Range("R9").Value = 0.09650924024640657
$ws.Range("S9").Value = 0.4353182751540041

# Row 10
$ws.Range("B10").Value = 0.09038382170862568
$ws.Range("D10").Value = 0.01981015270326042
$ws.Range("E10").Value = 0.001650846058605035
$ws.Range("F10").Value = 0.06397028477094512
$ws.Range("J10").Value = 0.116384647131655
$ws.Range("O10").Value = 0.01238134543953776
$ws.Range("Q10").Value = 0.2199752373091209
$ws.Range("R10").Value = 0.07924061081304168
$ws.Range("S10").Value = 0.3962030540652084

# Row 11
$ws.Range("G11").Value = 0.1239067055393586
$ws.Range("J11").Value = 0.09766763848396501
$ws.Range("K11").Value = 0.1865889212827988
$ws.Range("L11").Value = 0.5685131195335277
$ws.Range("S11").Value = 0.02332361516034985

# Row 12
$ws.Range("G12").Value = 0.7226277372262774
$ws.Range("J12").Value = 0.1800486618004866
$ws.Range("K12").Value = 0.0072992700729927
$ws.Range("L12").Value = 0.0389294403892944
$ws.Range("S12").Value = 0.05109489051094891

# Row 13
$ws.Range("G13").Value = 0.6588235294117647
$ws.Range("J13").Value = 0.2941176470588235
$ws.Range("S13").Value = 0.04705882352941176

# Row 15
$ws.Range("F15").Value = 0.0170940170940171
$ws.Range("H15").Value = 0.1773504273504274
$ws.Range("I15").Value = 0.07264957264957266
$ws.Range("J15").Value = 0.2948717948717949
$ws.Range("K15").Value = 0.07692307692307693
$ws.Range("M15").Value = 0.0170940170940171
$ws.Range("N15").Value = 0.004273504273504274
$ws.Range("O15").Value = 0.08547008547008547
$ws.Range("S15").Value = 0.2542735042735043

# Row 16
$ws.Range("F16").Value = 0.02686567164179104
$ws.Range("H16").Value = 0.217910447761194
$ws.Range("I16").Value = 0.09850746268656717
$ws.Range("J16").Value = 0.3671641791044776
$ws.Range("K16").Value = 0.1313432835820895
$ws.Range("M16").Value = 0.01194029850746269
$ws.Range("O16").Value = 0.03880597014925373
$ws.Range("S16").Value = 0.1074626865671642

# Row 17
$ws.Range("F17").Value = 0.01812366737739872
$ws.Range("H17").Value = 0.2025586353944563
$ws.Range("I17").Value = 0.09488272921108742
$ws.Range("J17").Value = 0.3678038379530917
$ws.Range("K17").Value = 0.1140724946695096
$ws.Range("M17").Value = 0.0138592750533049
$ws.Range("N17").Value = 0.003198294243070362
$ws.Range("O17").Value = 0.07462686567164178
$ws.Range("S17").Value = 0.1108742004264392

# Row 18
$ws.Range("F18").Value = 0.01030927835051546
$ws.Range("H18").Value = 0.1701030927835052
$ws.Range("I18").Value = 0.1108247422680412
$ws.Range("J18").Value = 0.3530927835051547
$ws.Range("K18").Value = 0.1262886597938144
$ws.Range("M18").Value = 0.03092783505154639
$ws.Range("O18").Value = 0.09020618556701031
$ws.Range("S18").Value = 0.1082474226804124

# Row 19
$ws.Range("F19").Value = 0.01505692251193537
$ws.Range("H19").Value = 0.2262210796915167
$ws.Range("I19").Value = 0.1053984575835476
$ws.Range("J19").Value = 0.3507161219243481
$ws.Range("K19").Value = 0.1160484759456482
$ws.Range("M19").Value = 0.01909658464928388
$ws.Range("N19").Value = 0.0007344840249724568
$ws.Range("O19").Value = 0.05728975394785164
$ws.Range("S19").Value = 0.1094381197208961
